$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 9, pushing the existing rows 9-12
# (TEST supplier / KLNMOP / GE rows) down to become rows 11-14.
$ws.Rows("9:10").Insert()

# Copy the formatting (borders/style) from row 8 into the two new rows
# so the new cells pick up the same style index (s="2") as the rest of
# the data rows instead of the default style.
$ws.Range("A8:I8").Copy()
$ws.Range("A9:I10").PasteSpecial(-4122)

# Populate the new Supplier row 9 ("TEST1 SUP")
$ws.Range("B9").Value = "TEST1 SUP"
$ws.Range("B10").Value = "TEST2 SUP"
$ws.Range("D9").Value = "OL"
$ws.Range("D10").Value = "OL"

$ws.Range("C9").Value = 535353
$ws.Range("E9").Value = "301-777"
$ws.Range("F9").Value = "Moscow"
$ws.Range("G9").Value = "Livs- 555"
$ws.Range("H9").Value = "Moscow@NoName.fr"
$ws.Range("I9").Value = $true

# Populate the new Supplier row 10 ("TEST2 SUP")
$ws.Range("C10").Value = 656565
$ws.Range("E10").Value = "101-777"
$ws.Range("F10").Value = "New York"
$ws.Range("G10").Value = "Topolowa- 555"
$ws.Range("H10").Value = "Moscow@NoName.fr"
$ws.Range("I10").Value = $true
